$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - sheet1
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 3403
$wsExpo.Range("F3").Value = 21
$wsExpo.Range("F5").Value = 1619
$wsExpo.Range("F6").Value = 63
$wsExpo.Range("F7").Value = 330

# Sheet "全部类型" (All types) - sheet4
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 3403
$wsAll.Range("F3").Value = 21
$wsAll.Range("F5").Value = 1619
$wsAll.Range("F6").Value = 63
$wsAll.Range("F8").Value = 330
